$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 4), mirroring the existing row 3 layout/format.
# Copy row 3's cells first so number formats / styles (date format on A,
# boolean type on B/G/I) carry over, then overwrite with the new values.
$ws.Range("A3:I3").Copy($ws.Range("A4:I4"))

$ws.Range("A4").Value = 42633.679085648146
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9974
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.22
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -0.52
$ws.Range("I4").Value = $false
